# Adds the "Live Templates" table (columns I/J) to the IntelliJ Shortcuts sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# ---- Header: I1:J1 merged, bold, left aligned ("Live Templates") ----
$ws.Range("I1").Value = "Live Templates"
$ws.Range("I1").Font.Bold = $true
$ws.Range("I1:J1").HorizontalAlignment = -4131
$ws.Range("I1:J1").Merge()

# ---- Body rows: shortcut (col I) / description (col J) ----
$rows = @(
    @{ Row = 2;  I = "sout";  J = "System.out" },
    @{ Row = 3;  I = "soutm"; J = "System.out method name" },
    @{ Row = 4;  I = "soutp"; J = "System.out paramaters" },
    @{ Row = 5;  I = "soutv"; J = "System.out value" },
    @{ Row = 6;  I = "souf";  J = "System.out.printf" },
    @{ Row = 8;  I = "psvm / main"; J = "public static void main" },
    @{ Row = 9;  I = "psf";   J = "public static final" },
    @{ Row = 10; I = "prsf";  J = "private static final" },
    @{ Row = 12; I = "iter";  J = "Iterator" },
    @{ Row = 13; I = "itco";  J = "Iterate collection" },
    @{ Row = 14; I = "itar";  J = "Iterate array" },
    @{ Row = 15; I = "itli";  J = "Iterate list" },
    @{ Row = 16; I = "fori";  J = "for loop" }
)

foreach ($r in $rows) {
    $cellI = $ws.Cells.Item($r.Row, 9)
    $cellJ = $ws.Cells.Item($r.Row, 10)
    $cellI.Value = $r.I
    $cellJ.Value = $r.J
    $cellI.Font.ColorIndex = -4105
    $cellJ.Font.ColorIndex = -4105
}

# Last row (18) uses the same highlighted style as the other "last entry in a
# group" cells on this sheet (light-blue fill + box border), so copy that
# existing format (from G19, the label-style half of the pair) onto both
# I18 and J18.
$ws.Range("I18").Value = "call"
$ws.Range("J18").Value = "Surround with method call"
$ws.Range("G19").Copy() | Out-Null
$ws.Range("I18").PasteSpecial(-4122) | Out-Null
$ws.Range("G19").Copy() | Out-Null
$ws.Range("J18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---- Column J width (bestfit-like, ~24 chars) ----
$ws.Columns.Item(10).ColumnWidth = 23.15

# ---- Selection moves to I8 (matches the author's last-edited cell) ----
$ws.Range("I8").Select()

Write-Output "Live Templates table added"
